$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Key"
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"

# Data row
$ws.Range("A2").Value = "RegisterWithoutLastName"
$ws.Range("B2").Value = "Iliya"

# Update selection to match target state
$ws.Range("C2").Select()
